$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.831729054450989
$ws.Range("B1").Value = 2.3274986743927
$ws.Range("C1").Value = 2.469321727752686
$ws.Range("D1").Value = 2.933242559432983
$ws.Range("E1").Value = 1.998267769813538
